# Remove the last column ("Employee") from the user-groups/roles table.
#
# The table starts life with 7 columns:
#   1560 | 1985 | 1219 | 1219 | 1219 | 1219 | 1219
# (blank, Function, System admin, Reception, Management, Supervisor, Employee)
#
# The edit drops the trailing "Employee" column (header cell + every data
# cell in the 16 rows), drops the matching <w:gridCol>, and shrinks the
# table's overall preferred width by the width of that column
# (9640 - 1219 = 8421 twips).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$lastColIndex = $t.Columns.Count
$t.Columns.Item($lastColIndex).Delete()

# Table.PreferredWidth is expressed in points (dxa = points * 20), so the
# desired 8421 dxa table width becomes 421.05 points.
$t.PreferredWidth = 421.05
